$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the date header cells (C2:C6)
$ws.Range("C2").Value = "07/08/2019"
$ws.Range("C3").Value = "07/09/2019"
$ws.Range("C4").Value = "07/10/2019"
$ws.Range("C5").Value = "07/11/2019"
$ws.Range("C6").Value = "07/12/2019"

# Update the remarks text cells (D2:D6)
$ws.Range("D6").Value = "1.`tCreated defect #7907 for the production issue(Access Denied on CreateProfile.aspx page). Tried to reproduce the issue in different ways, but not able to reproduce it here. Please refer the mail sent for details.`n2.`tCreated defect #7908 as per yesterday’s observations shared with Sonny and assigned it to him in PMall admin.`n3.`tWorked on creation of Life Cycle for a Project/Feature/Defect in PMall admin. Shared the initial draft in separate email. Please review and suggest.`n4.`tTried to reproduce the other production issue reported yesterday for iPad, but not able to reproduce that too. Please refer the email shared for more details. Also, attached the same here for reference.`n5.`tAttended daily status meeting to discuss project updates, tasks , observations and requirements.`n6.`tUpdated testing status report for June release. Please find attached sheet with updated Ticket-Wise Test Summary, Test Case Execution Summary, Test cases and Defects created in June release till date."
$ws.Range("D5").Value = "1.`tRetested defect #7746 again on dev environment which is failing retest on dev as well for issue 2. Assigned it to you for review. Please review and suggest.`n2.`tCreated defects #7900, #7901, #7902 in PMall admin as per Sonny’s response on yesterday’s observations.`n3.`tRetested defects #7815 and #7817. Please respond on query shared regarding defect #7817.`n4.`tVerified “Re-order This Item” functionality from Order Details page on mobile site as part of #6791. Shared observations with Sonny for review.`n5.`tAttended daily status meeting to discuss project updates, tasks , observations and requirements.`n6.`tUpdated testing status report for June release. Please find attached sheet with updated Ticket-Wise Test Summary, Test Case Execution Summary, Test cases and Defects created in June release till date. Also shared detailed status in a separate email."
$ws.Range("D4").Value = "1.`tPerformed regression test on live environment on both mobile and desktop site. No issues found. Shared status in a separate email.`n2.`tRetested defect #7746 on dev environment which is failing retest. Assigned it back to Sonny for review. Will update the same as per our discussion in today’s stand-up call.`n3.`tVerified 3 new pages as part of #6791 on dev environment. Shared observations with Sonny and you on email. Please review and let me know if defects need to be created against them.`n4.`tReviewed errors appearing in Bug snag tool. Need to discuss the same to identify the priority and where to raise them.`n5.`tAttended daily status meeting to discuss project updates, tasks , observations and requirements.`n6.`tWill reschedule daily stand-up call to 11:15 AM(Chicago time), yesterday onwards.`n7.`tUpdated testing status report for June release. Please find attached sheet with updated Ticket-Wise Test Summary, Test Case Execution Summary, Test cases and Defects created in June release till date. Also shared detailed status in a separate email."
$ws.Range("D3").Value = "1.`tRetested defects assigned to me in PMall admin. Shared detailed status in a separate email, attached here for reference.`n2.`tCreated defects #7890 and #7891 based on discussion on yesterday’s observations and assigned them to Sonny.`n3.`tWorked on verification of more pages as per Sonny’s confirmation on scope of ticket #6791 - Mobile round corner block refresh. Will share observations if any once complete.`n4.`tUpdated testing status report for June release. Please find attached sheet with updated Ticket-Wise Test Summary, Test Case Execution Summary, Test cases and Defects created in June release till date. Also shared detailed status in a separate email."
$ws.Range("D2").Value = "1.`tPerformed ad-hoc testing on different tickets deployed on Dev environment as part of June Release.`n2.`tUpdated testing status report for June release. Please find attached sheet with updated Ticket-Wise Test Summary, Test Case Execution Summary, Test cases and Defects created in June release till date. Also shared detailed status in a separate email.`n3.`tAttended daily status meeting to discuss project updates, tasks , observations and requirements."

# Update row heights
$ws.Rows.Item(2).RowHeight = 105
$ws.Rows.Item(3).RowHeight = 150
$ws.Rows.Item(4).RowHeight = 225
$ws.Rows.Item(5).RowHeight = 195
$ws.Rows.Item(6).RowHeight = 225
